# Update "想去人数" (number of people interested) counts on the three
# sheets that share data with the 展览 sheet, matching the regenerated
# gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 108
$ws.Range("F6").Value = 5272
$ws.Range("F7").Value = 443
$ws.Range("F9").Value = 910
$ws.Range("F11").Value = 73
$ws.Range("F12").Value = 28
$ws.Range("F13").Value = 566
$ws.Range("F14").Value = 14
$ws.Range("F15").Value = 17
$ws.Range("F17").Value = 1759
$ws.Range("F19").Value = 828
$ws.Range("F22").Value = 305
$ws.Range("F23").Value = 510
$ws.Range("F24").Value = 130
$ws.Range("F28").Value = 2575
$ws.Range("F30").Value = 98
$ws.Range("F34").Value = 253
$ws.Range("F40").Value = 639
$ws.Range("F42").Value = 46

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 232

# --- Sheet "全部类型" (all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 232
$ws.Range("F4").Value = 108
$ws.Range("F7").Value = 5272
$ws.Range("F8").Value = 443
$ws.Range("F12").Value = 910
$ws.Range("F16").Value = 73
$ws.Range("F17").Value = 28
$ws.Range("F18").Value = 566
$ws.Range("F19").Value = 14
$ws.Range("F20").Value = 17
$ws.Range("F23").Value = 1759
$ws.Range("F25").Value = 828
$ws.Range("F27").Value = 305
$ws.Range("F29").Value = 510
$ws.Range("F30").Value = 130
$ws.Range("F33").Value = 2575
$ws.Range("F35").Value = 98
$ws.Range("F38").Value = 253
$ws.Range("F44").Value = 46
